$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-25: fix up individual "D" (column E on sheet) imputed/missing values ---
# E2: RM 2 -> was missing, now -7.2
$ws.Range("E2").Value = -7.2
# E6: RM 21 -> was -5.7, now missing
$ws.Range("E6").Formula = '=""'
# E12: RM 81 -> was missing, now -5.3
$ws.Range("E12").Value = -5.3
# E14: RM 90 -> was -5.4, now missing
$ws.Range("E14").Formula = '=""'
# E20: RM 134 -> was missing, now -7.2
$ws.Range("E20").Value = -7.2
# E21: RM 135 -> was missing, now -8.7
$ws.Range("E21").Value = -8.7
# E23: RM 140 -> was -7, now missing
$ws.Range("E23").Formula = '=""'
# E24: RM 142a -> was -8.1, now missing
$ws.Range("E24").Formula = '=""'

# --- Rows 26-35: "RM 232" and "SC 92" rows removed; remaining SC rows shift up,
#     and a handful of D/E values change along the way. Overwrite rows 26-33
#     with their final contents directly, then drop the two now-unused trailing rows. ---

# Row 26 (was SC 5 @ row 27)
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

# Row 27 (was SC 101 @ row 29)
$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Formula = '=""'
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

# Row 28 (was SC 105 @ row 30)
$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

# Row 29 (was SC 119 @ row 31)
$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

# Row 30 (was SC 120 @ row 32)
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

# Row 31 (was SC 132 @ row 33)
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

# Row 32 (was SC 193 @ row 34)
$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Formula = '=""'
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

# Row 33 (was SC 232 @ row 35)
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

# Drop the two now-redundant trailing rows (this happens at the very bottom of the
# used range so nothing below needs to shift).
$ws.Rows("34:35").Delete()
